# Updates cryptos price/volume columns (D, E) for rows 2-51
# to the refreshed values captured in this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.212.78' },
    @{ Cell = 'E2'; Value = '  -2.38%  ' },
    @{ Cell = 'D3'; Value = '1.822.77' },
    @{ Cell = 'E3'; Value = '  -1.94%  ' },
    @{ Cell = 'E4'; Value = '  -1.23%  ' },
    @{ Cell = 'D5'; Value = '314.49' },
    @{ Cell = 'E5'; Value = '  -1.90%  ' },
    @{ Cell = 'D6'; Value = '1.005' },
    @{ Cell = 'E6'; Value = '  -1.04%  ' },
    @{ Cell = 'D7'; Value = '0.4265' },
    @{ Cell = 'E7'; Value = '  -1.67%  ' },
    @{ Cell = 'D8'; Value = '0.3683' },
    @{ Cell = 'E8'; Value = '  -2.92%  ' },
    @{ Cell = 'D9'; Value = '0.07248' },
    @{ Cell = 'E9'; Value = '  -2.53%  ' },
    @{ Cell = 'D10'; Value = '0.8621' },
    @{ Cell = 'E10'; Value = '  -2.83%  ' },
    @{ Cell = 'D11'; Value = '21.01' },
    @{ Cell = 'E11'; Value = '  -3.24%  ' },
    @{ Cell = 'D12'; Value = '1.829.32' },
    @{ Cell = 'E12'; Value = '  -1.72%  ' },
    @{ Cell = 'D13'; Value = '6.722' },
    @{ Cell = 'E13'; Value = '  -1.10%  ' },
    @{ Cell = 'D14'; Value = '0.07096' },
    @{ Cell = 'E14'; Value = '  -0.72%  ' },
    @{ Cell = 'D15'; Value = '5.311' },
    @{ Cell = 'E15'; Value = '  -3.49%  ' },
    @{ Cell = 'D16'; Value = '88.15' },
    @{ Cell = 'E16'; Value = '  -0.19%  ' },
    @{ Cell = 'D17'; Value = '1.007' },
    @{ Cell = 'E17'; Value = '  -1.34%  ' },
    @{ Cell = 'D18'; Value = '0.000008869' },
    @{ Cell = 'E18'; Value = '  -2.06%  ' },
    @{ Cell = 'D19'; Value = '1.006' },
    @{ Cell = 'E19'; Value = '  -1.02%  ' },
    @{ Cell = 'D20'; Value = '15.06' },
    @{ Cell = 'E20'; Value = '  -3.10%  ' },
    @{ Cell = 'D21'; Value = '27.235.06' },
    @{ Cell = 'E21'; Value = '  -2.32%  ' },
    @{ Cell = 'D22'; Value = '5.137' },
    @{ Cell = 'E22'; Value = '  -2.76%  ' },
    @{ Cell = 'D23'; Value = '10.86' },
    @{ Cell = 'E23'; Value = '  -3.13%  ' },
    @{ Cell = 'D24'; Value = '2.045.47' },
    @{ Cell = 'E24'; Value = '  -2.08%  ' },
    @{ Cell = 'D25'; Value = '2.002' },
    @{ Cell = 'E25'; Value = '  -1.41%  ' },
    @{ Cell = 'D26'; Value = '153.31' },
    @{ Cell = 'E26'; Value = '  -2.35%  ' },
    @{ Cell = 'D27'; Value = '18.35' },
    @{ Cell = 'E27'; Value = '  -1.69%  ' },
    @{ Cell = 'D28'; Value = '2.141' },
    @{ Cell = 'E28'; Value = '  +6.18%  ' },
    @{ Cell = 'D29'; Value = '5.230' },
    @{ Cell = 'E29'; Value = '  -3.73%  ' },
    @{ Cell = 'D30'; Value = '116.34' },
    @{ Cell = 'E30'; Value = '  -3.36%  ' },
    @{ Cell = 'D31'; Value = '0.08895' },
    @{ Cell = 'E31'; Value = '  -1.26%  ' },
    @{ Cell = 'D32'; Value = '1.195' },
    @{ Cell = 'E32'; Value = '  -3.79%  ' },
    @{ Cell = 'D33'; Value = '0.7545' },
    @{ Cell = 'E33'; Value = '  -2.70%  ' },
    @{ Cell = 'D34'; Value = '4.429' },
    @{ Cell = 'E34'; Value = '  -3.59%  ' },
    @{ Cell = 'D35'; Value = '2.807' },
    @{ Cell = 'E35'; Value = '  -4.93%  ' },
    @{ Cell = 'D36'; Value = '1.006' },
    @{ Cell = 'E36'; Value = '  -1.03%  ' },
    @{ Cell = 'D37'; Value = '1.111' },
    @{ Cell = 'E37'; Value = '  -2.66%  ' },
    @{ Cell = 'D38'; Value = '0.01972' },
    @{ Cell = 'E38'; Value = '  -0.20%  ' },
    @{ Cell = 'D39'; Value = '0.05265' },
    @{ Cell = 'E39'; Value = '  -0.93%  ' },
    @{ Cell = 'D40'; Value = '7.158' },
    @{ Cell = 'E40'; Value = '  +2.04%  ' },
    @{ Cell = 'D41'; Value = '2.865' },
    @{ Cell = 'E41'; Value = '  -0.58%  ' },
    @{ Cell = 'D42'; Value = '0.1695' },
    @{ Cell = 'E42'; Value = '  +0.70%  ' },
    @{ Cell = 'D43'; Value = '0.5039' },
    @{ Cell = 'E43'; Value = '  -3.33%  ' },
    @{ Cell = 'D44'; Value = '8.632' },
    @{ Cell = 'E44'; Value = '  -1.90%  ' },
    @{ Cell = 'D45'; Value = '10.60' },
    @{ Cell = 'E45'; Value = '  -1.94%  ' },
    @{ Cell = 'D46'; Value = '107.29' },
    @{ Cell = 'E46'; Value = '  -3.17%  ' },
    @{ Cell = 'D47'; Value = '0.4739' },
    @{ Cell = 'E47'; Value = '  -0.51%  ' },
    @{ Cell = 'E48'; Value = '  -1.07%  ' },
    @{ Cell = 'D49'; Value = '0.06372' },
    @{ Cell = 'E49'; Value = '  -1.94%  ' },
    @{ Cell = 'D50'; Value = '1.659' },
    @{ Cell = 'E50'; Value = '  -3.29%  ' },
    @{ Cell = 'D51'; Value = '1.812' },
    @{ Cell = 'E51'; Value = '  -4.43%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text interpretation so numeric-looking strings (e.g. "21.01")
    # are not silently coerced into floating point numbers by Excel,
    # then restore the default "Normal" style so no stray format sticks.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

